$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Metric (column C) and Metric - Sort (column D) values for rows 55-83 ---
# This reflects the insertion of new "Residential Aged Care ..." weekly metrics and
# the corresponding re-sequencing of every metric below it in the Aged Care block.

$ws.Range("C55").Value = "# Aged Care Active Resident Cases (Weekly)"
$ws.Range("D55").Value = 550
$ws.Range("C56").Value = "# Aged Care Active Resident Cases (Weekly) per 1M"
$ws.Range("D56").Value = 560
$ws.Range("C57").Value = "% Aged Care Active Resident Cases (Weekly) Change"
$ws.Range("D57").Value = 570
$ws.Range("C58").Value = "# Aged Care Staff Cases"
$ws.Range("D58").Value = 580
$ws.Range("C59").Value = "# Aged Care Staff Cases (7-day avg)"
$ws.Range("D59").Value = 590
$ws.Range("C60").Value = "# Aged Care Staff Cases (7-day avg) per 1M"
$ws.Range("D60").Value = 600
$ws.Range("C61").Value = "% Aged Care Staff Cases Weekly Change"
$ws.Range("D61").Value = 610
$ws.Range("C62").Value = "# Aged Care Active Staff Cases (Weekly)"
$ws.Range("D62").Value = 640
$ws.Range("C63").Value = "# Aged Care Active Staff Cases (Weekly) per 1M"
$ws.Range("D63").Value = 650
$ws.Range("C64").Value = "% Aged Care Active Staff Cases (Weekly) Change"
$ws.Range("D64").Value = 660
$ws.Range("C65").Value = "# Aged Care Active Outbreaks"
$ws.Range("D65").Value = 670
$ws.Range("C66").Value = "# Aged Care Active Outbreaks (7-day avg)"
$ws.Range("D66").Value = 680
$ws.Range("C67").Value = "# Aged Care Active Outbreaks (7-day avg) per 1M"
$ws.Range("D67").Value = 690
$ws.Range("C68").Value = "% Aged Care Active Outbreaks Weekly Change"
$ws.Range("D68").Value = 700
$ws.Range("C69").Value = "# Aged Care Outbreaks Facilities (Weekly)"
$ws.Range("D69").Value = 710
$ws.Range("C70").Value = "# Aged Care Outbreaks Facilities (Weekly) per 1M"
$ws.Range("D70").Value = 720
$ws.Range("C71").Value = "% Aged Care Outbreaks Facilities (Weekly) Change"
$ws.Range("D71").Value = 730
$ws.Range("C72").Value = "# Residential Aged Care Resident Cases (Weekly)"
$ws.Range("D72").Value = 740
$ws.Range("C73").Value = "# Residential Aged Care Staff Cases (Weekly)"
$ws.Range("D73").Value = 750
$ws.Range("C74").Value = "# Residential Aged Care Staff Cases (Weekly) per 1M"
$ws.Range("D74").Value = 760
$ws.Range("C75").Value = "# Aged Care Molnupiravir Prescriptions (Daily)"
$ws.Range("D75").Value = 770
$ws.Range("C76").Value = "# Aged Care Molnupiravir Prescriptions (Daily) per 1M"
$ws.Range("D76").Value = 780
$ws.Range("C77").Value = "% Aged Care Molnupiravir Prescriptions (Daily) per Case"
$ws.Range("D77").Value = 790
$ws.Range("C78").Value = "# Aged Care Paxlovid Prescriptions (Daily)"
$ws.Range("D78").Value = 800
$ws.Range("C79").Value = "# Aged Care Paxlovid Prescriptions (Daily) per 1M"
$ws.Range("D79").Value = 810
$ws.Range("C80").Value = "% Aged Care Paxlovid Prescriptions (Daily) per Case"
$ws.Range("D80").Value = 820
$ws.Range("C81").Value = "# Monthly PBS Scripts"
$ws.Range("D81").Value = 830
$ws.Range("C82").Value = "# Monthly PBS Scripts per 1M"
$ws.Range("D82").Value = 840
$ws.Range("C83").Value = "% Monthly PBS Scripts Change"
$ws.Range("D83").Value = 850

# --- Style (s="2") changes on column C: moves down three rows within the block ---
# Rows 65-67 lose the s="2" formatting; rows 69-71 gain it.
$ws.Range("C68").Copy()
$ws.Range("C65:C67").PasteSpecial(-4122)
$ws.Range("C72").Copy()
$ws.Range("C69:C71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the active selection / view to match the saved workbook state ---
$ws.Activate() | Out-Null
$ws.Range("A61").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D79").Select() | Out-Null
